# Weekly refresh of the "Hortaliza, Vega Monumental Concepción - Zanahoria" sheet.
# A new week's observation (Primera/Segunda) is prepended at rows 83-84, which
# pushes every existing record (previously rows 83-136) down by two rows,
# growing the used range from A1:R136 to A1:R138.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the first data row that changes; Excel
# shifts rows 83:136 down to 85:138, which reproduces the rest of the diff
# (every existing row's content simply slides down by two rows) without
# having to rewrite ~50 untouched rows by hand.
$ws.Rows("83:84").Insert()

# New "Primera" observation.
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44489
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100114013
$ws.Range("G83").Value = "Zanahoria"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 600
$ws.Range("K83").Value = 7000
$ws.Range("L83").Value = 7500
$ws.Range("M83").Value = 7250
$ws.Range("N83").Value = "`$/saco 20 kilos"
$ws.Range("O83").Value = "Región de Ñuble"
$ws.Range("P83").Value = 362
$ws.Range("Q83").Value = 20
$ws.Range("R83").Value = "Hortaliza"

# New "Segunda" observation.
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 44489
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100114013
$ws.Range("G84").Value = "Zanahoria"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 300
$ws.Range("K84").Value = 6500
$ws.Range("L84").Value = 6500
$ws.Range("M84").Value = 6500
$ws.Range("N84").Value = "`$/saco 20 kilos"
$ws.Range("O84").Value = "Región de Ñuble"
$ws.Range("P84").Value = 325
$ws.Range("Q84").Value = 20
$ws.Range("R84").Value = "Hortaliza"
